# Update summary tables with newest airtoxics nata data
$wb = $excel.ActiveWorkbook

# --- Sheet "Means" ---
$wsMeans = $wb.Worksheets.Item("Means")

# Row 9: Total Cancer Risk (per million)
$wsMeans.Range("B9").Value = 26
$wsMeans.Range("C9").Value = 34
$wsMeans.Range("D9").Value = 30
$wsMeans.Range("E9").Value = 35
$wsMeans.Range("F9").Value = 34
$wsMeans.Range("G9").Value = 35

# Row 10: Total Respiratory (hazard quotient)
$wsMeans.Range("B10").Value = 0.31
$wsMeans.Range("C10").Value = 0.39
$wsMeans.Range("D10").Value = 0.4
$wsMeans.Range("E10").Value = 0.42
$wsMeans.Range("F10").Value = 0.41
$wsMeans.Range("G10").Value = 0.41

# --- Sheet "Standard Deviations" ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# Row 9: Total Cancer Risk (per million)
$wsSD.Range("B9").Value = 8.3
$wsSD.Range("C9").Value = 9.3
$wsSD.Range("E9").Value = 5.3
$wsSD.Range("F9").Value = 5
$wsSD.Range("G9").Value = 5

# Row 10: Total Respiratory (hazard quotient)
$wsSD.Range("B10").Value = 0.11
$wsSD.Range("E10").Value = 0.046
$wsSD.Range("F10").Value = 0.036
$wsSD.Range("G10").Value = 0.029
